$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Amend two old polls (rows 4 and 5) and fill in the CAN (H) column for a
# newly added poll (row 6 here corresponds to the NDP line of the poll that
# already had most of its other province columns filled), plus refresh the
# two sample-size rows (nw/nu) and the recalculated vote-sum check.
# The "CAN" (H) column values change style from the explicit Normal-with-font
# xf to the workbook's bare default xf (applyFont=false), matching how these
# particular cells were (re)entered, so clear formatting to General first.
$ws.Range("H4:H9").NumberFormat = "General"

$ws.Range("H4").Value = 35
$ws.Range("H5").Value = 29
$ws.Range("H6").Value = 23
$ws.Range("H7").Value = 7
$ws.Range("H8").Value = 5
$ws.Range("H9").Value = 2

$ws.Range("H10").Value = 1223
$ws.Range("H11").Value = 1248

$null = $ws.Range("H12").Select()
